# Fix the csv missing fields for Aquitanian mapping file
#
# Populates the previously-empty "review" (I) and "dob" (J) columns for
# every data row (2-14) with the placeholder value "[]", and moves the
# active selection to I15:J15 (next empty row in those columns), matching
# the author's worksheet-view state after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 14

for ($row = 2; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 9).Value  = "[]"   # column I - review
    $ws.Cells.Item($row, 10).Value = "[]"   # column J - dob
}

# Scroll the view roughly to where the new data lives and select the next
# empty I:J cells, mirroring the saved sheetView/selection state.
$win = $excel.ActiveWindow
$win.ScrollRow = 11
$win.ScrollColumn = 5

$ws.Range("I15:J15").Select()
